# update northern ireland data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Refresh the Northern Ireland row (row 6) with the latest published figures
$ws.Range("B6").Value = 44008
$ws.Range("C6").Value = 824
$ws.Range("D6").Value = 427
$ws.Range("E6").Value = 343
$ws.Range("F6").Value = 48
$ws.Range("G6").Value = 8

# The old "Republic of Korea" row (row 22) is no longer needed - remove it
$ws.Rows.Item(22).Delete()

# Restore the on-open view to the top of the sheet
$null = $ws.Range("F10").Select()
